# Updates crypto price/volume data (and swaps Cosmos/InjectiveProtocol rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.277.55"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "3.584.87"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  -0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "605.32"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.03%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "147.70"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("D7").Value = "3.583.39"
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").Value = "4.193.28"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("E14").Value = "  -0.29%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "29.52"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("D16").Value = "3.583.56"
$ws.Range("E16").Value = "  +1.04%  "
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").Value = "66.351.39"
$ws.Range("E18").Value = "  +0.29%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.06"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.47%  "
$ws.Range("E20").Value = "  +2.25%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.83"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.20%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "423.16"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.54%  "
$ws.Range("E23").Value = "  +0.30%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "78.40"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  +3.29%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.20"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +4.82%  "
$ws.Range("E28").Value = "  +2.91%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "3.582.54"
$ws.Range("E31").Value = "  +0.89%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.156"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +3.58%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "25.03"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("E38").Value = "  -2.71%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "174.49"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.47%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0853"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.04%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "5.18"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  -1.04%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "45.87"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("E44").Value = "  -3.62%  "
$ws.Range("E45").Value = "  +0.09%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.53"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +5.57%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "23.51"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.60%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "24.21"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.86%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "7.14"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("E50").Value = "  -5.06%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.953"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.58%  "
